# Applies the daily-scrape update: adds 11 new opportunity rows (rows 2-12)
# and widens columns B:H to fit the newly scraped content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
# Excel pads ColumnWidth (chars) by ~0.8333 (5/6) char when it stores the
# <col width=.../> value, so subtract that padding to land on the exact
# target stored widths (B=55, C=86, D=65, E=10, F=16, G=15, H=34).
$colPadding = 0.8333333333333333
$ws.Columns.Item(2).ColumnWidth = 55 - $colPadding
$ws.Columns.Item(3).ColumnWidth = 86 - $colPadding
$ws.Columns.Item(4).ColumnWidth = 65 - $colPadding
$ws.Columns.Item(5).ColumnWidth = 10 - $colPadding
$ws.Columns.Item(6).ColumnWidth = 16 - $colPadding
$ws.Columns.Item(7).ColumnWidth = 15 - $colPadding
$ws.Columns.Item(8).ColumnWidth = 34 - $colPadding

# --- New data rows -------------------------------------------------------
# OPPORTUNITY ID values are numeric-looking but must stay text, so each A
# cell is formatted as Text right before its value is written (matches the
# scraped source data, which keeps the IDs as strings).

# Row 2: Sales Lead
$ws.Range('A2').NumberFormat = '@'
$ws.Range('A2').Value = '1328826'
$ws.Range('B2').Value = 'https://aiesec.org/opportunity/global-talent/1328826'
$ws.Range('C2').Value = 'Sales Lead'
$ws.Range('D2').Value = 'Nairobi, Kenya'
$ws.Range('E2').Value = 'No'
$ws.Range('F2').Value = '2 applicants'
$ws.Range('G2').Value = '3 - 6 Months'
$ws.Range('H2').Value = 'YOUR APPS LIMITED'

# Row 3: Power Electronics Internship involving Inverter Experimentat
$ws.Range('A3').NumberFormat = '@'
$ws.Range('A3').Value = '1328731'
$ws.Range('B3').Value = 'https://aiesec.org/opportunity/global-talent/1328731'
$ws.Range('C3').Value = 'Power Electronics Internship involving Inverter Experimentation and PSIM Simulation'
$ws.Range('D3').Value = '日本、兵庫県神戸市'
$ws.Range('E3').Value = 'No'
$ws.Range('F3').Value = '1 applicant'
$ws.Range('G3').Value = '9 - 12 Weeks'
$ws.Range('H3').Value = 'Sohatsu Systems Laboratory Inc.'

# Row 4: Site Engineer
$ws.Range('A4').NumberFormat = '@'
$ws.Range('A4').Value = '1328691'
$ws.Range('B4').Value = 'https://aiesec.org/opportunity/global-talent/1328691'
$ws.Range('C4').Value = 'Site Engineer'
$ws.Range('D4').Value = 'Novi Sad, Serbia'
$ws.Range('E4').Value = 'No'
$ws.Range('F4').Value = '0 applicants'
$ws.Range('G4').Value = '9 - 12 Weeks'
$ws.Range('H4').Value = 'Orbis Company'

# Row 5: Accelerate Romania | Marketing & Social Media Manager
$ws.Range('A5').NumberFormat = '@'
$ws.Range('A5').Value = '1328456'
$ws.Range('B5').Value = 'https://aiesec.org/opportunity/global-talent/1328456'
$ws.Range('C5').Value = 'Accelerate Romania | Marketing & Social Media Manager'
$ws.Range('D5').Value = 'Arad, Romania'
$ws.Range('E5').Value = 'No'
$ws.Range('F5').Value = '1 applicant'
$ws.Range('G5').Value = '9 - 12 Weeks'
$ws.Range('H5').Value = 'BNB'

# Row 6: Guest Relations Officer Intern
$ws.Range('A6').NumberFormat = '@'
$ws.Range('A6').Value = '1328021'
$ws.Range('B6').Value = 'https://aiesec.org/opportunity/global-talent/1328021'
$ws.Range('C6').Value = 'Guest Relations Officer Intern'
$ws.Range('D6').Value = 'Nugegoda, Sri Lanka'
$ws.Range('E6').Value = 'No'
$ws.Range('F6').Value = '0 applicants'
$ws.Range('G6').Value = '3 - 6 Months'
$ws.Range('H6').Value = 'The Barn By Starbeans in Ella'

# Row 7: Marketing Specialist
$ws.Range('A7').NumberFormat = '@'
$ws.Range('A7').Value = '1327239'
$ws.Range('B7').Value = 'https://aiesec.org/opportunity/global-talent/1327239'
$ws.Range('C7').Value = 'Marketing Specialist'
$ws.Range('D7').Value = 'Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt'
$ws.Range('E7').Value = 'No'
$ws.Range('F7').Value = '16 applicants'
$ws.Range('G7').Value = '3 - 6 Months'
$ws.Range('H7').Value = 'Karcel'

# Row 8: Videographer & Video Editor
$ws.Range('A8').NumberFormat = '@'
$ws.Range('A8').Value = '1327236'
$ws.Range('B8').Value = 'https://aiesec.org/opportunity/global-talent/1327236'
$ws.Range('C8').Value = 'Videographer & Video Editor'
$ws.Range('D8').Value = 'Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt'
$ws.Range('E8').Value = 'No'
$ws.Range('F8').Value = '1 applicant'
$ws.Range('G8').Value = '3 - 6 Months'
$ws.Range('H8').Value = 'Karcel'

# Row 9: Content Creator (Storyteller & social media Maven)
$ws.Range('A9').NumberFormat = '@'
$ws.Range('A9').Value = '1327232'
$ws.Range('B9').Value = 'https://aiesec.org/opportunity/global-talent/1327232'
$ws.Range('C9').Value = 'Content Creator (Storyteller & social media Maven)'
$ws.Range('D9').Value = 'Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt'
$ws.Range('E9').Value = 'No'
$ws.Range('F9').Value = '4 applicants'
$ws.Range('G9').Value = '3 - 6 Months'
$ws.Range('H9').Value = 'Karcel'

# Row 10: Graphic Designer
$ws.Range('A10').NumberFormat = '@'
$ws.Range('A10').Value = '1327124'
$ws.Range('B10').Value = 'https://aiesec.org/opportunity/global-talent/1327124'
$ws.Range('C10').Value = 'Graphic Designer'
$ws.Range('D10').Value = 'Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt'
$ws.Range('E10').Value = 'No'
$ws.Range('F10').Value = '2 applicants'
$ws.Range('G10').Value = '3 - 6 Months'
$ws.Range('H10').Value = 'Karcel'

# Row 11: Social Media Planner & Content Creator
$ws.Range('A11').NumberFormat = '@'
$ws.Range('A11').Value = '1325318'
$ws.Range('B11').Value = 'https://aiesec.org/opportunity/global-talent/1325318'
$ws.Range('C11').Value = 'Social Media Planner & Content Creator'
$ws.Range('D11').Value = 'القاهرة، محافظة القاهرة‬، مصر'
$ws.Range('E11').Value = 'No'
$ws.Range('F11').Value = '19 applicants'
$ws.Range('G11').Value = '9 - 12 Weeks'
$ws.Range('H11').Value = 'The Circle Care'

# Row 12: Digital Marketing Executive
$ws.Range('A12').NumberFormat = '@'
$ws.Range('A12').Value = '1321400'
$ws.Range('B12').Value = 'https://aiesec.org/opportunity/global-talent/1321400'
$ws.Range('C12').Value = 'Digital Marketing Executive'
$ws.Range('D12').Value = 'Cairo, Cairo Governorate, Egypt'
$ws.Range('E12').Value = 'No'
$ws.Range('F12').Value = '34 applicants'
$ws.Range('G12').Value = '3 - 6 Months'
$ws.Range('H12').Value = 'Silverkey Technologies Egypt'
